$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: update title (D9) and link (E9)
$ws.Range("D9").Value = "상아탑과 현실이 동떨어진 나라"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/ivory-tower-real-world/#utm_source=rss&utm_medium=rss&utm_campaign=ivory-tower-real-world"

# Row 28: update title (D28) and link (E28)
$ws.Range("D28").Value = "Let's do MuJoCo - 3.1 Kinematics (기구학을 공부하는데 도움 되었던 자료)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/177"
